$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AA2").Value = 'maa://21246 (91.23), maa://36684 (98.63), ***maa://22731 (6.67)'
$ws.Range("K3").Value = '*maa://22880 (69.68), maa://20276 (82.86), *maa://22749 (66.67)'
$ws.Range("W3").Value = 'maa://27396 (85.02), maa://27484 (95.74), maa://27480 (82.35)'
$ws.Range("AA3").Value = 'maa://24390 (96.15)'
$ws.Range("S4").Value = 'maa://32509 (98.78), maa://22754 (91.67), maa://27295 (80.39), *maa://21746 (55.81), *maa://31008 (78.05)'
$ws.Range("W4").Value = '**maa://32495 (47.54), ***maa://31785 (19.47), ***maa://36683 (26.67)'
$ws.Range("AA6").Value = 'maa://22739 (91.67)'
$ws.Range("W7").Value = 'maa://22399 (94.66), *maa://22758 (70.37)'
$ws.Range("AE7").Value = '*maa://26191 (68.49), *maa://36671 (72.73), maa://42530 (100.0)'
$ws.Range("O8").Value = 'maa://32931 (88.89), *maa://21916 (60.34), maa://23252 (92.31), **maa://22759 (45.45), maa://37496 (100.0)'
$ws.Range("W8").Value = 'maa://21411 (96.04)'
$ws.Range("W9").Value = 'maa://26223 (97.0)'
$ws.Range("AA9").Value = 'maa://28711 (88.1), ***maa://22740 (5.88), **maa://27377 (46.15), ***maa://25174 (20.0), **maa://39938 (43.75), maa://40166 (100.0)'
$ws.Range("W10").Value = 'maa://22301 (97.42), maa://22726 (100.0)'
$ws.Range("W11").Value = 'maa://36713 (97.86)'
$ws.Range("G12").Value = 'maa://21867 (90.07)'
$ws.Range("G13").Value = '*maa://21248 (75.48), **maa://22728 (47.62)'
$ws.Range("C15").Value = '*maa://22743 (76.61), maa://22734 (83.33), *maa://30808 (64.29), ***maa://36048 (12.9)'
$ws.Range("G15").Value = 'maa://24304 (88.11), maa://21478 (91.18)'
$ws.Range("AE15").Value = 'maa://21364 (80.68), *maa://22766 (73.0), *maa://36666 (78.46)'
$ws.Range("AE16").Value = '*maa://23911 (61.96), maa://27755 (91.89)'
$ws.Range("G17").Value = 'maa://22430 (88.14), maa://39599 (80.95)'
$ws.Range("S17").Value = '***maa://42324 (28.57)'
$ws.Range("C18").Value = 'maa://24570 (96.65)'
$ws.Range("AA19").Value = '*maa://30709 (60.9), *maa://36668 (52.17)'
$ws.Range("G20").Value = 'maa://22864 (88.55)'
$ws.Range("K23").Value = 'maa://39756 (92.59), maa://39875 (95.83)'
$ws.Range("W23").Value = '*maa://28503 (63.93)'
$ws.Range("C24").Value = 'maa://24368 (80.42)'
$ws.Range("W24").Value = 'maa://23504 (92.92), maa://29988 (86.34), **maa://22892 (40.14), *maa://25141 (77.05), maa://36663 (80.7), ***maa://22815 (23.08)'
$ws.Range("C25").Value = 'maa://29753 (95.18)'
$ws.Range("G25").Value = '*maa://29063 (76.12), *maa://25311 (74.19), ***maa://22725 (4.84)'
$ws.Range("AA26").Value = '*maa://42235 (76.19)'
$ws.Range("AE27").Value = 'maa://24023 (96.83)'
$ws.Range("W28").Value = 'maa://39929 (86.83), ***maa://39723 (14.71), maa://41749 (81.25)'
$ws.Range("AE28").Value = 'maa://36660 (93.87), *maa://36701 (64.0)'
$ws.Range("K29").Value = 'maa://28432 (93.52), *maa://28440 (72.84), maa://31400 (100.0), *maa://28650 (66.67)'
$ws.Range("AE29").Value = '*maa://24080 (69.02), ***maa://34960 (8.7)'
$ws.Range("S32").Value = 'maa://41108 (89.47), maa://41238 (94.59)'
$ws.Range("S34").Value = 'maa://24526 (93.16)'
$ws.Range("K35").Value = 'maa://41296 (98.15)'
$ws.Range("AE38").Value = 'maa://36697 (84.73)'
$ws.Range("O39").Value = 'maa://24709 (92.23)'
$ws.Range("O40").Value = 'maa://23278 (95.89), maa://21386 (95.63), maa://36664 (90.24)'
$ws.Range("G44").Value = 'maa://29768 (97.54), maa://27728 (96.0)'
$ws.Range("G46").Value = 'maa://35931 (92.61)'
$ws.Range("G47").Value = 'maa://27410 (95.82), maa://29661 (97.64), maa://28038 (84.62)'
$ws.Range("G53").Value = 'maa://32534 (93.26), **maa://32434 (36.36)'
